# Weekly refresh of the "Pepino ensalada" consolidated sheet:
# a new observation is inserted as row 107 (pushing the existing
# rows 107-153 down to 108-154), and the fresh row is populated
# with this period's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107, shifting rows 107:153 down to 108:154.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new weekly record.
$ws.Cells.Item(107, 1).Value = 11
$ws.Cells.Item(107, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(107, 3).Value = "Bíobío"
$ws.Cells.Item(107, 4).Value = 44813
$ws.Cells.Item(107, 5).Value = 8
$ws.Cells.Item(107, 6).Value = 100112043
$ws.Cells.Item(107, 7).Value = "Pepino ensalada"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 100
$ws.Cells.Item(107, 11).Value = 24000
$ws.Cells.Item(107, 12).Value = 25000
$ws.Cells.Item(107, 13).Value = 24500
$ws.Cells.Item(107, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(107, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(107, 16).Value = 408
$ws.Cells.Item(107, 17).Value = 60
$ws.Cells.Item(107, 18).Value = "Hortaliza"
